$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "60,0"
$ws.Range("B1").Value = "190,0"
$ws.Range("C1").Value = "7,0"
$ws.Range("D1").Value = "IMC: 16,62"
$ws.Range("E1").Value = "Bajo peso"

# Column widths re-flow to fit the new (longer/shorter) cell contents, same as
# Excel's bestFit auto-sizing would do when the sheet is regenerated.
$ws.Columns.Item(1).ColumnWidth = 4.9296875
$ws.Columns.Item(2).ColumnWidth = 6.046875
$ws.Columns.Item(3).ColumnWidth = 3.81640625
$ws.Columns.Item(4).ColumnWidth = 10.7421875
$ws.Columns.Item(5).ColumnWidth = 9.7421875
